$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Drop the hidden "version/cover" bookmark that wrapped the title
#    run (it is a hidden bookmark - name starts with "_" - so it is
#    not enumerated by Bookmarks, but it can still be reached by name).
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_30j0zll")) {
    $d.Bookmarks("_30j0zll").Delete()
}

# ------------------------------------------------------------------
# 2. Remove the whole version/cover table (the "Versão / Data /
#    Descrição" table that float-anchors over the first paragraph).
# ------------------------------------------------------------------
while ($d.Tables.Count -gt 0) {
    $d.Tables(1).Delete()
}

# ------------------------------------------------------------------
# 3. Turn the old "<HAIR2U>" subtitle paragraph into the new
#    underlined "Representação Sistemica" heading paragraph.
# ------------------------------------------------------------------
$p1 = $d.Paragraphs(1)

# Drop the "Subtitle" style -> plain (Normal) paragraph.
$p1.Style = $d.Styles("Normal")

# Paragraph formatting: widowControl off, single (240 auto) spacing.
$p1.Format.WidowControl = 0
$p1.Format.LineSpacingRule = 0

# Replace the run text.
$p1.Range.Text = "Representação Sistemica"

# Underline only the visible text run, not the paragraph mark.
$p1 = $d.Paragraphs(1)
$textRange = $p1.Range
[void]$textRange.MoveEnd(1, -1)
$textRange.Font.Underline = 1

# ------------------------------------------------------------------
# 4. The paragraph that used to sit right after the table (empty,
#    but carrying widowControl/spacing/underline formatting) is now
#    redundant - its formatting has been folded into paragraph 1.
# ------------------------------------------------------------------
$p2 = $d.Paragraphs(2)
if ($p2.Range.Text -eq [char]13) {
    $p2.Range.Delete()
}
